$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.562.29"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.462.30"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -1.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.63"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.14"
$ws.Range("E6").Value = "  +2.92%  "
$ws.Range("E7").Value = "  +3.57%  "
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.511"
$ws.Range("E9").Value = "  +7.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.47"
$ws.Range("E10").Value = "  +2.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0793"
$ws.Range("E11").Value = "  +4.56%  "
$ws.Range("E12").Value = "  +1.42%  "
$ws.Range("D13").Value = "2.842.59"
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.86"
$ws.Range("E14").Value = "  +2.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.78"
$ws.Range("E15").Value = "  +6.47%  "
$ws.Range("D16").Value = "2.497.04"
$ws.Range("E16").Value = "  +3.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.772"
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("D18").Value = "41.558.67"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.50"
$ws.Range("E19").Value = "  +6.29%  "
$ws.Range("D20").Value = "0.0₃0937"
$ws.Range("E20").Value = "  +5.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.82"
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("E22").Value = "  +5.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.80"
$ws.Range("E23").Value = "  +2.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.71"
$ws.Range("E24").Value = "  +2.67%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.90"
$ws.Range("E26").Value = "  +4.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.25"
$ws.Range("E27").Value = "  +4.29%  "
$ws.Range("E28").Value = "  +2.25%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.62"
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.11"
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.36"
$ws.Range("E31").Value = "  +3.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.43"
$ws.Range("E32").Value = "  +3.68%  "
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0756"
$ws.Range("E34").Value = "  +3.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.17"
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.41"
$ws.Range("E36").Value = "  -5.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.85"
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("E38").Value = "  +3.62%  "
$ws.Range("E39").Value = "  +5.06%  "
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("E41").Value = "  +2.96%  "
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("D43").Value = "1.957.24"
$ws.Range("E43").Value = "  +1.99%  "
$ws.Range("E44").Value = "  +3.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.72"
$ws.Range("E45").Value = "  -8.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.90"
$ws.Range("E46").Value = "  +2.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.96"
$ws.Range("E47").Value = "  +5.76%  "
$ws.Range("D48").Value = "2.701.95"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.39"
$ws.Range("E49").Value = "  +3.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.80"
$ws.Range("E50").Value = "  +4.46%  "
$ws.Range("E51").Value = "  +0.67%  "
